$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

Set-TextValue $ws.Range("D2") "43.035.78"
$ws.Range("E2").Value = "  -4.98%  "

Set-TextValue $ws.Range("D3") "2.226.11"
$ws.Range("E3").Value = "  -5.70%  "

$ws.Range("E4").Value = "  -0.02%  "

Set-TextValue $ws.Range("D5") "318.09"
$ws.Range("E5").Value = "  +2.46%  "

Set-TextValue $ws.Range("D6") "98.97"
$ws.Range("E6").Value = "  -9.26%  "

Set-TextValue $ws.Range("D7") "0.580"
$ws.Range("E7").Value = "  -7.58%  "

$ws.Range("E8").Value = "  -0.02%  "

Set-TextValue $ws.Range("D9") "0.565"
$ws.Range("E9").Value = "  -8.20%  "

Set-TextValue $ws.Range("D10") "36.51"
$ws.Range("E10").Value = "  -11.39%  "

Set-TextValue $ws.Range("D11") "54.36"
$ws.Range("E11").Value = "  -1.93%  "

Set-TextValue $ws.Range("D12") "0.0825"
$ws.Range("E12").Value = "  -9.99%  "

Set-TextValue $ws.Range("D13") "7.72"
$ws.Range("E13").Value = "  -8.76%  "

$ws.Range("E14").Value = "  -4.07%  "

$ws.Range("B15").Value = "Polygon"
$ws.Range("C15").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
Set-TextValue $ws.Range("D15") "0.866"
$ws.Range("E15").Value = "  -12.20%  "

$ws.Range("B16").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C16").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
Set-TextValue $ws.Range("D16") "2.564.22"
$ws.Range("E16").Value = "  -5.90%  "

Set-TextValue $ws.Range("D17") "14.04"
$ws.Range("E17").Value = "  -8.59%  "

Set-TextValue $ws.Range("D18") "2.213.29"
$ws.Range("E18").Value = "  -6.31%  "

Set-TextValue $ws.Range("D19") "42.808.62"
$ws.Range("E19").Value = "  -5.45%  "

Set-TextValue $ws.Range("D20") "14.46"
$ws.Range("E20").Value = "  +4.59%  "

Set-TextValue $ws.Range("D21") "0.0₃0963"
$ws.Range("E21").Value = "  -9.42%  "

Set-TextValue $ws.Range("D22") "6.46"
$ws.Range("E22").Value = "  -11.93%  "

Set-TextValue $ws.Range("D23") "65.15"
$ws.Range("E23").Value = "  -11.19%  "

$ws.Range("E24").Value = "  -9.80%  "

Set-TextValue $ws.Range("D25") "235.89"
$ws.Range("E25").Value = "  -9.07%  "

$ws.Range("E26").Value = "  -8.96%  "

Set-TextValue $ws.Range("D27") "1.00"
$ws.Range("E27").Value = "  +0.28%  "

Set-TextValue $ws.Range("D28") "10.16"
$ws.Range("E28").Value = "  -8.89%  "

$ws.Range("E29").Value = "  -7.10%  "

$ws.Range("E30").Value = "  -14.33%  "

Set-TextValue $ws.Range("D31") "0.0883"
$ws.Range("E31").Value = "  -8.65%  "

Set-TextValue $ws.Range("D32") "20.46"
$ws.Range("E32").Value = "  -8.31%  "

Set-TextValue $ws.Range("D33") "157.33"
$ws.Range("E33").Value = "  -6.95%  "

Set-TextValue $ws.Range("D34") "33.75"
$ws.Range("E34").Value = "  -10.97%  "

$ws.Range("E35").Value = "  -5.83%  "

Set-TextValue $ws.Range("D36") "3.35"
$ws.Range("E36").Value = "  +12.95%  "

Set-TextValue $ws.Range("D37") "2.03"
$ws.Range("E37").Value = "  +16.49%  "

$ws.Range("E38").Value = "  -6.36%  "

Set-TextValue $ws.Range("D39") "4.48"
$ws.Range("E39").Value = "  -7.33%  "

Set-TextValue $ws.Range("D40") "0.104"
$ws.Range("E40").Value = "  -10.84%  "

Set-TextValue $ws.Range("D41") "3.64"
$ws.Range("E41").Value = "  -7.46%  "

Set-TextValue $ws.Range("D42") "0.0324"
$ws.Range("E42").Value = "  -9.19%  "

Set-TextValue $ws.Range("D43") "1.852.09"
$ws.Range("E43").Value = "  +10.60%  "

$ws.Range("E44").Value = "  +0.06%  "

Set-TextValue $ws.Range("D45") "12.16"
$ws.Range("E45").Value = "  -5.81%  "

Set-TextValue $ws.Range("D46") "88.14"
$ws.Range("E46").Value = "  -11.13%  "

Set-TextValue $ws.Range("D47") "5.49"
$ws.Range("E47").Value = "  -0.39%  "

Set-TextValue $ws.Range("D48") "78.57"
$ws.Range("E48").Value = "  -5.14%  "

Set-TextValue $ws.Range("D49") "0.206"
$ws.Range("E49").Value = "  -11.33%  "

Set-TextValue $ws.Range("D50") "60.48"
$ws.Range("E50").Value = "  -13.36%  "

Set-TextValue $ws.Range("D51") "8.65"
$ws.Range("E51").Value = "  -5.25%  "
